$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new blank row right after the header row, pushing all existing
# data rows (old row 2 .. old row 102) down by one (new row 3 .. new row 103).
$ws.Rows("2:2").Insert()

# The inserted row picks up formatting from the header row by default;
# reset it to the plain "Normal" style used by the rest of the data rows,
# then re-apply the date number format on column D (same format the other
# "Fecha" cells in this column use).
$ws.Range("A2:R2").Style = "Normal"
$ws.Range("D2").NumberFormat = "YYYY-MM-DD HH:MM:SS"

# Populate the new record (week of 2023-01-30, Región del Maule, Primera).
$ws.Range("A2").Value = 10
$ws.Range("B2").Value = "Vega Modelo de Temuco"
$ws.Range("C2").Value = "La Araucanía"
$ws.Range("D2").Value = 44956
$ws.Range("E2").Value = 9
$ws.Range("F2").Value = 100112030
$ws.Range("G2").Value = "Poroto granado"
$ws.Range("H2").Value = "Sin especificar"
$ws.Range("I2").Value = "Primera"
$ws.Range("J2").Value = 50
$ws.Range("K2").Value = 50000
$ws.Range("L2").Value = 50000
$ws.Range("M2").Value = 50000
$ws.Range("N2").Value = "$/saco 25 kilos"
$ws.Range("O2").Value = "Región del Maule"
$ws.Range("P2").Value = 2000
$ws.Range("Q2").Value = 25
$ws.Range("R2").Value = "Hortaliza"
